$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $style = $c.Style
    $c.Value = "'" + $val
    $c.Style = $style
}

Set-TextValue "D2" "65.491.29"
Set-TextValue "E2" "  +2.35%  "
Set-TextValue "D3" "3.190.88"
Set-TextValue "E3" "  +5.65%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "569.55"
Set-TextValue "E5" "  +3.97%  "
Set-TextValue "D6" "148.72"
Set-TextValue "E6" "  +9.11%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "D8" "3.183.87"
Set-TextValue "E8" "  +5.88%  "
Set-TextValue "D9" "0.506"
Set-TextValue "E9" "  +5.10%  "
Set-TextValue "D10" "6.90"
Set-TextValue "E10" "  +9.09%  "
Set-TextValue "E11" "  +6.23%  "
Set-TextValue "D12" "0.480"
Set-TextValue "E12" "  +6.63%  "
Set-TextValue "D13" "38.01"
Set-TextValue "E13" "  +7.89%  "
Set-TextValue "D14" "0.0000229"
Set-TextValue "E14" "  +5.90%  "
Set-TextValue "D15" "3.702.26"
Set-TextValue "E15" "  +5.77%  "
Set-TextValue "D16" "65.602.34"
Set-TextValue "E16" "  +2.44%  "
Set-TextValue "B17" "TRON"
Set-TextValue "C17" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D17" "0.114"
Set-TextValue "E17" "  +2.94%  "
Set-TextValue "B18" "WrappedEther"
Set-TextValue "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "3.189.03"
Set-TextValue "E18" "  +5.56%  "
Set-TextValue "B19" "BitcoinCash"
Set-TextValue "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "530.54"
Set-TextValue "E19" "  +11.83%  "
Set-TextValue "D20" "7.04"
Set-TextValue "E20" "  +8.91%  "
Set-TextValue "D21" "14.41"
Set-TextValue "E21" "  +7.68%  "
Set-TextValue "D22" "0.734"
Set-TextValue "E22" "  +9.32%  "
Set-TextValue "D23" "7.67"
Set-TextValue "E23" "  +10.34%  "
Set-TextValue "D24" "13.24"
Set-TextValue "E24" "  +8.11%  "
Set-TextValue "D25" "80.37"
Set-TextValue "E25" "  +3.99%  "
Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  -0.08%  "
Set-TextValue "D27" "9.28"
Set-TextValue "E27" "  +22.42%  "
Set-TextValue "E28" "  +8.48%  "
Set-TextValue "E29" "  +8.86%  "
Set-TextValue "D30" "27.20"
Set-TextValue "E30" "  +7.16%  "
Set-TextValue "D31" "0.999"
Set-TextValue "E31" "  -0.14%  "
Set-TextValue "E32" "  +4.02%  "
Set-TextValue "D33" "1.16"
Set-TextValue "E33" "  +5.39%  "
Set-TextValue "D34" "555.16"
Set-TextValue "E34" "  -4.69%  "
Set-TextValue "B35" "NEARProtocol"
Set-TextValue "C35" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D35" "5.56"
Set-TextValue "E35" "  +5.97%  "
Set-TextValue "B36" "Filecoin"
Set-TextValue "C36" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D36" "6.26"
Set-TextValue "E36" "  +9.24%  "
Set-TextValue "D37" "54.55"
Set-TextValue "E37" "  +5.98%  "
Set-TextValue "E38" "  +8.12%  "
Set-TextValue "D39" "0.0846"
Set-TextValue "E39" "  +8.35%  "
Set-TextValue "E40" "  +7.45%  "
Set-TextValue "D41" "3.193.58"
Set-TextValue "E41" "  +10.31%  "
Set-TextValue "E42" "  +5.86%  "
Set-TextValue "D43" "8.51"
Set-TextValue "E43" "  +5.39%  "
Set-TextValue "D44" "0.277"
Set-TextValue "E44" "  +17.05%  "
Set-TextValue "D45" "2.31"
Set-TextValue "E45" "  +13.64%  "
Set-TextValue "D46" "26.52"
Set-TextValue "E46" "  +8.86%  "
Set-TextValue "D47" "1.00"
Set-TextValue "E47" "  +0.04%  "
Set-TextValue "B48" "Monero"
Set-TextValue "C48" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D48" "123.73"
Set-TextValue "E48" "  +5.21%  "
Set-TextValue "B49" "PEPE"
Set-TextValue "C49" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D49" "0.0₃0541"
Set-TextValue "E49" "  +4.20%  "
Set-TextValue "E50" "  +4.48%  "
Set-TextValue "D51" "2.17"
